$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 878.5357
$ws.Range("I98").Value = 892.55554
$ws.Range("K98").Value = 892.55554
$ws.Range("M98").Value = 605.44446
$ws.Range("H122").Value = 878.5357
$ws.Range("I122").Value = 892.55554
$ws.Range("K122").Value = 2677.66662
$ws.Range("M122").Value = -227.66662
$ws.Range("H124").Value = 38780
$ws.Range("J124").Value = 38780
$ws.Range("L124").Value = 38780
$ws.Range("N124").Value = -48600
$ws.Range("H125").Value = 1545.1666
$ws.Range("I125").Value = 1721.2
$ws.Range("J125").Value = 665
$ws.Range("K125").Value = 15490.8
$ws.Range("L125").Value = 5985
$ws.Range("M125").Value = -13030.8
$ws.Range("N125").Value = -10905
$ws.Range("H126").Value = 45000
$ws.Range("J126").Value = 45000
$ws.Range("L126").Value = 45000
$ws.Range("N126").Value = -54880
$ws.Range("H127").Value = 2404.2727
$ws.Range("I127").Value = 497
$ws.Range("J127").Value = 2595
$ws.Range("K127").Value = 1491
$ws.Range("L127").Value = 7785
$ws.Range("M127").Value = 3469
$ws.Range("N127").Value = -17705
$ws.Range("H138").Value = 4506771.5
$ws.Range("J138").Value = 3590.4482
$ws.Range("L138").Value = 10771.3446
$ws.Range("N138").Value = -21051.3446

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2109.311
$ws.Range("I61").Value = 2164.7437
$ws.Range("K61").Value = 2164.7437
$ws.Range("M61").Value = -1952.7437
$ws.Range("H136").Value = 2109.311
$ws.Range("I136").Value = 2164.7437
$ws.Range("K136").Value = 6494.2311
$ws.Range("M136").Value = -3944.2311

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 44007.43
$ws.Range("J21").Value = 44007.43
$ws.Range("L21").Value = 44007.43
$ws.Range("N21").Value = -44479.43
$ws.Range("H118").Value = 7891.837
$ws.Range("J118").Value = 7891.837
$ws.Range("L118").Value = 7891.837
$ws.Range("N118").Value = -11205.837
$ws.Range("H130").Value = 39995
$ws.Range("J130").Value = 39995
$ws.Range("L130").Value = 39995
$ws.Range("N130").Value = -50035
$ws.Range("H134").Value = 3396.238
$ws.Range("I134").Value = 2308.842
$ws.Range("J134").Value = 5049.08
$ws.Range("K134").Value = 6926.526
$ws.Range("L134").Value = 15147.24
$ws.Range("M134").Value = -4391.526
$ws.Range("N134").Value = -20217.24

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2678.1228
$ws.Range("I31").Value = 2236.4707
$ws.Range("J31").Value = 3331
$ws.Range("K31").Value = 2236.4707
$ws.Range("L31").Value = 3331
$ws.Range("M31").Value = -1941.4707
$ws.Range("N31").Value = -3921
$ws.Range("H34").Value = 2678.1228
$ws.Range("I34").Value = 2236.4707
$ws.Range("J34").Value = 3331
$ws.Range("K34").Value = 2236.4707
$ws.Range("L34").Value = 3331
$ws.Range("M34").Value = -2034.4707
$ws.Range("N34").Value = -3735
$ws.Range("H75").Value = 42513.332
$ws.Range("J75").Value = 42513.332
$ws.Range("L75").Value = 42513.332
$ws.Range("N75").Value = -44509.332
$ws.Range("H78").Value = 42513.332
$ws.Range("J78").Value = 42513.332
$ws.Range("L78").Value = 127539.996
$ws.Range("N78").Value = -137523.996
$ws.Range("H130").Value = 70000
$ws.Range("J130").Value = 70000
$ws.Range("L130").Value = 70000
$ws.Range("N130").Value = -80040

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 38462304
$ws.Range("I113").Value = 200000580
$ws.Range("J113").Value = 814
$ws.Range("K113").Value = 600001740
$ws.Range("L113").Value = 2442
$ws.Range("M113").Value = -599999570
$ws.Range("N113").Value = -6782
$ws.Range("H122").Value = 1243.8667
$ws.Range("I122").Value = 1111.9412
$ws.Range("J122").Value = 1416.3846
$ws.Range("K122").Value = 10007.4708
$ws.Range("L122").Value = 12747.4614
$ws.Range("M122").Value = -7557.470799999999
$ws.Range("N122").Value = -17647.4614

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3109.276
$ws.Range("I102").Value = 4210.294
$ws.Range("J102").Value = 1549.5
$ws.Range("K102").Value = 4210.294
$ws.Range("L102").Value = 1549.5
$ws.Range("M102").Value = -2588.294
$ws.Range("N102").Value = -4793.5
$ws.Range("H113").Value = 144492.14
$ws.Range("I113").Value = 251261.25
$ws.Range("J113").Value = 2133.3333
$ws.Range("K113").Value = 251261.25
$ws.Range("L113").Value = 2133.3333
$ws.Range("M113").Value = -249091.25
$ws.Range("N113").Value = -6473.3333
$ws.Range("H122").Value = 2408.75
$ws.Range("I122").Value = 2454.1
$ws.Range("K122").Value = 7362.299999999999
$ws.Range("M122").Value = -4912.299999999999
$ws.Range("H124").Value = 60780
$ws.Range("J124").Value = 60780
$ws.Range("L124").Value = 60780
$ws.Range("N124").Value = -70600
$ws.Range("H126").Value = 3054.3333
$ws.Range("I126").Value = 2177.6667
$ws.Range("J126").Value = 3580.3333
$ws.Range("K126").Value = 6533.000100000001
$ws.Range("L126").Value = 10740.9999
$ws.Range("M126").Value = -4063.000100000001
$ws.Range("N126").Value = -15680.9999
$ws.Range("H128").Value = 50779.75
$ws.Range("J128").Value = 50779.75
$ws.Range("L128").Value = 50779.75
$ws.Range("N128").Value = -60739.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7547.0586
$ws.Range("I7").Value = 18666.666
$ws.Range("J7").Value = 5164.2856
$ws.Range("K7").Value = 18666.666
$ws.Range("L7").Value = 5164.2856
$ws.Range("M7").Value = -18554.666
$ws.Range("N7").Value = -5388.2856
$ws.Range("H36").Value = 41715
$ws.Range("J36").Value = 41715
$ws.Range("L36").Value = 41715
$ws.Range("N36").Value = -42839
$ws.Range("H40").Value = 9029
$ws.Range("I40").Value = 18233.334
$ws.Range("J40").Value = 2125.75
$ws.Range("K40").Value = 18233.334
$ws.Range("L40").Value = 2125.75
$ws.Range("M40").Value = -18097.334
$ws.Range("N40").Value = -2397.75
$ws.Range("H122").Value = 9057.058999999999
$ws.Range("I122").Value = 9510
$ws.Range("J122").Value = 8740
$ws.Range("K122").Value = 28530
$ws.Range("L122").Value = 26220
$ws.Range("M122").Value = -26080
$ws.Range("N122").Value = -31120
$ws.Range("H126").Value = 7547.0586
$ws.Range("I126").Value = 18666.666
$ws.Range("J126").Value = 5164.2856
$ws.Range("K126").Value = 55999.99800000001
$ws.Range("L126").Value = 15492.8568
$ws.Range("M126").Value = -53529.99800000001
$ws.Range("N126").Value = -20432.8568
$ws.Range("H127").Value = 48308.08
$ws.Range("J127").Value = 48308.08
$ws.Range("L127").Value = 48308.08
$ws.Range("N127").Value = -58228.08

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 38682.5
$ws.Range("J75").Value = 38682.5
$ws.Range("L75").Value = 38682.5
$ws.Range("N75").Value = -40554.5
$ws.Range("H78").Value = 38682.5
$ws.Range("J78").Value = 38682.5
$ws.Range("L78").Value = 116047.5
$ws.Range("N78").Value = -125407.5
$ws.Range("H122").Value = 2348.7715
$ws.Range("I122").Value = 2400.1428
$ws.Range("J122").Value = 2143.2856
$ws.Range("K122").Value = 7200.428400000001
$ws.Range("L122").Value = 6429.8568
$ws.Range("M122").Value = -4750.428400000001
$ws.Range("N122").Value = -11329.8568
$ws.Range("H126").Value = 1680.7059
$ws.Range("I126").Value = 1397.7241
$ws.Range("J126").Value = 3322
$ws.Range("K126").Value = 4193.1723
$ws.Range("L126").Value = 9966
$ws.Range("M126").Value = -1723.1723
$ws.Range("N126").Value = -14906
$ws.Range("H127").Value = 24304
$ws.Range("J127").Value = 24304
$ws.Range("L127").Value = 24304
$ws.Range("N127").Value = -34224
$ws.Range("H131").Value = 79795
$ws.Range("J131").Value = 79795
$ws.Range("L131").Value = 79795
$ws.Range("N131").Value = -89875
$ws.Range("H132").Value = 1382.4108
$ws.Range("I132").Value = 1200.1177
$ws.Range("K132").Value = 3600.3531
$ws.Range("M132").Value = -1070.3531
